$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 103, shifting existing rows 103:153 down to 104:154
$ws.Rows("103:103").Insert()

# Populate the new row 103 with the latest price-report entry
$ws.Range("A103").Value2 = 10
$ws.Range("B103").Value2 = "Vega Modelo de Temuco"
$ws.Range("C103").Value2 = "La Araucanía"
$ws.Range("D103").Value2 = 44489
$ws.Range("E103").Value2 = 9
$ws.Range("F103").Value2 = "Fruta"
$ws.Range("G103").Value2 = 100102
$ws.Range("H103").Value2 = "Cítricos"
$ws.Range("I103").Value2 = 100102006
$ws.Range("J103").Value2 = "Pomelo"
$ws.Range("K103").Value2 = "Start Ruby"
$ws.Range("L103").Value2 = "Especial"
$ws.Range("M103").Value2 = 70
$ws.Range("N103").Value2 = 15000
$ws.Range("O103").Value2 = 15000
$ws.Range("P103").Value2 = 15000
$ws.Range("Q103").Value2 = '$/bandeja 15 kilos granel'
$ws.Range("R103").Value2 = "Región de O'Higgins"
$ws.Range("S103").Value2 = 1000
$ws.Range("T103").Value2 = 15
